# commons_structure.xlsx - classFields sheet update
#
# The "domain.Order" and "domain.Order$OrderBuilder" field rows get
# re-sorted into a new field order, and the "domain.OrderStatus" enum
# constant rows plus two "domain.Topics" rows get reordered as well.
# (Field types / names below are the final values shown after the edit;
# target cell already holding the correct text is simply a no-op write.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# domain.Order fields (rows 2-8)
$ws.Range("B2").Value  = 'productId'
$ws.Range("D2").Value  = 'java.lang.Long'
$ws.Range("B3").Value  = 'source'
$ws.Range("D3").Value  = 'domain.OrderSource'
$ws.Range("B4").Value  = 'status'
$ws.Range("D4").Value  = 'domain.OrderStatus'
$ws.Range("B5").Value  = 'customerId'
$ws.Range("D5").Value  = 'java.lang.Long'
$ws.Range("B6").Value  = 'price'
$ws.Range("D6").Value  = 'int'
$ws.Range("B7").Value  = 'productCount'
$ws.Range("D7").Value  = 'int'
$ws.Range("B8").Value  = 'id'
$ws.Range("D8").Value  = 'java.lang.Long'

# domain.OrderStatus enum constants (rows 9-15)
$ws.Range("B9").Value  = 'ACCEPT'
$ws.Range("B10").Value = '$VALUES'
$ws.Range("B11").Value = 'ROLLBACK'
$ws.Range("B12").Value = 'REJECT'
$ws.Range("B13").Value = 'REJECTED'
$ws.Range("B14").Value = 'CONFIRMED'
$ws.Range("B15").Value = 'NEW'

# domain.Order$OrderBuilder fields (rows 18-24)
$ws.Range("B18").Value = 'customerId'
$ws.Range("D18").Value = 'java.lang.Long'
$ws.Range("B19").Value = 'productId'
$ws.Range("D19").Value = 'java.lang.Long'
$ws.Range("B20").Value = 'price'
$ws.Range("D20").Value = 'int'
$ws.Range("B21").Value = 'productCount'
$ws.Range("D21").Value = 'int'
$ws.Range("B22").Value = 'status'
$ws.Range("D22").Value = 'domain.OrderStatus'
$ws.Range("B23").Value = 'source'
$ws.Range("D23").Value = 'domain.OrderSource'
$ws.Range("B24").Value = 'id'
$ws.Range("D24").Value = 'java.lang.Long'

# domain.Topics constants (rows 29 & 31 swap)
$ws.Range("B29").Value = 'STOCK'
$ws.Range("B31").Value = 'ORDERS'
